# Weekly update: insert two new "Lechuga" (Escarola) price records at the
# top of the historical block (rows 979-980), pushing all existing rows
# 979-1046 down by two (to 981-1048).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 979 (each Insert() shifts rows 979.. down by one).
$ws.Rows.Item(979).Insert()
$ws.Rows.Item(979).Insert()

# Columns A, B, C, E, F, G, R are constant for every data row in this sheet.
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$categoriaId = 100112033
$categoria = "Lechuga"
$clasificacion = "Hortaliza"

# New row 979: Escarola / Primera, fecha 2023-12-05 (serial 45265)
$ws.Cells.Item(979, 1).Value = $mercadoId
$ws.Cells.Item(979, 2).Value = $mercado
$ws.Cells.Item(979, 3).Value = $region
$ws.Cells.Item(979, 4).Value = 45265
$ws.Cells.Item(979, 5).Value = $codreg
$ws.Cells.Item(979, 6).Value = $categoriaId
$ws.Cells.Item(979, 7).Value = $categoria
$ws.Cells.Item(979, 8).Value = "Escarola"
$ws.Cells.Item(979, 9).Value = "Primera"
$ws.Cells.Item(979, 10).Value = 350
$ws.Cells.Item(979, 11).Value = 16000
$ws.Cells.Item(979, 12).Value = 16000
$ws.Cells.Item(979, 13).Value = 16000
$ws.Cells.Item(979, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(979, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(979, 16).Value = 1067
$ws.Cells.Item(979, 17).Value = 15
$ws.Cells.Item(979, 18).Value = $clasificacion

# New row 980: Escarola / Segunda, fecha 2023-12-05 (serial 45265)
$ws.Cells.Item(980, 1).Value = $mercadoId
$ws.Cells.Item(980, 2).Value = $mercado
$ws.Cells.Item(980, 3).Value = $region
$ws.Cells.Item(980, 4).Value = 45265
$ws.Cells.Item(980, 5).Value = $codreg
$ws.Cells.Item(980, 6).Value = $categoriaId
$ws.Cells.Item(980, 7).Value = $categoria
$ws.Cells.Item(980, 8).Value = "Escarola"
$ws.Cells.Item(980, 9).Value = "Segunda"
$ws.Cells.Item(980, 10).Value = 350
$ws.Cells.Item(980, 11).Value = 14000
$ws.Cells.Item(980, 12).Value = 14000
$ws.Cells.Item(980, 13).Value = 14000
$ws.Cells.Item(980, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(980, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(980, 16).Value = 778
$ws.Cells.Item(980, 17).Value = 18
$ws.Cells.Item(980, 18).Value = $clasificacion
